# feat: add 2022-Q3 data
#
# 1) Insert a brand-new "2022-Q3" sheet (placed right after "总计", right
#    before "2022-Q2") holding the fund-holding detail rows for the new
#    quarter.
# 2) Update the "总计" (totals) sheet: shift all existing quarter summary
#    rows down by one and insert the 2022-Q3 summary as the new first row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" detail sheet.
#
# Duplicating the existing "2022-Q2" sheet (instead of Worksheets.Add)
# keeps the exact same cell styles (bold/centered headers + border on the
# header row and on the index column) without minting new style records.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The duplicate still has all 26 of 2022-Q2's fund rows (rows 2-27).
# 2022-Q3 only has 3 funds, so drop rows 5-27, leaving rows 2-4 free to be
# overwritten below.
$extra = $q3.Range($q3.Rows.Item(5), $q3.Rows.Item(27))
$extra.Delete()

# Fund code / name / size / position columns are stored as text in this
# workbook (even though several of them look numeric), so force the
# "text" number format before writing them, then clear the format again
# afterwards so the cells end up with no explicit style (matching the
# other detail sheets) while keeping the text type.
$q3.Range("B2:G4").NumberFormat = "@"

$q3.Range("B2").Value = "512480"
$q3.Range("C2").Value = "国联安半导体ETF"
$q3.Range("D2").Value = "122.28"
$q3.Range("E2").Value = "98.75"
$q3.Range("F2").Value = "2.94"
$q3.Range("G2").Value = "3.5950"
$q3.Range("H2").Value = 8

$q3.Range("B3").Value = "008327"
$q3.Range("C3").Value = "西藏东财中证通信技术主题指数C"
$q3.Range("D3").Value = "1.30"
$q3.Range("E3").Value = "94.83"
$q3.Range("F3").Value = "3.76"
$q3.Range("G3").Value = "0.0489"
$q3.Range("H3").Value = 8

$q3.Range("B4").Value = "008326"
$q3.Range("C4").Value = "西藏东财中证通信技术主题指数A"
$q3.Range("D4").Value = "0.39"
$q3.Range("E4").Value = "94.83"
$q3.Range("F4").Value = "3.76"
$q3.Range("G4").Value = "0.0147"
$q3.Range("H4").Value = 8

$q3.Range("B2:G4").ClearFormats()

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet.
#
# Column A is just the running 0-based row index, so it doesn't need to
# move - only the quarter label / count / market-value columns (B:D) shift
# down by one row. Shift from the bottom up so we don't clobber a row
# before it has been copied down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $srcRow = $total.Range("B" + $r + ":D" + $r)
    $dstRow = $total.Range("B" + ($r + 1) + ":D" + ($r + 1))
    $srcRow.Copy($dstRow)
}

# Extend the index column with the new last row (copy A8's style onto A9
# first so it keeps the bold/centered look without creating a new style).
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)
$total.Range("A9").Value = 7

# Fill in the new 2022-Q3 summary row at the top of the data.
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 3.66

# Restore the original active sheet/selection ("总计" was the active tab
# before this edit).
$total.Activate()
[void]$total.Range("A1").Select()

Write-Host "2022-Q3 sheet and 总计 summary updated"
